$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New values for rows 2-5, columns B-P (A is a label, unchanged)
$data = @{
    2 = @(4.11, 100, 11.72, 85.41, 55.95, 95.98, 88.95999999999999, 56.14, 48.59, 21.62, 12.92, 64.91, 18.63, 57.44, 98.37)
    3 = @(2.06, 100, 4.6, 91.18000000000001, 50.56, 98.41, 94, 52.81, 29.55, 26.74, 6.47, 38.84, 14.26, 53.1, 93.81)
    4 = @(5.71, 100, 5.71, 100, 94.29000000000001, 100, 97.14, 14.29, 82.86, 65.70999999999999, 40, 82.86, 74.29000000000001, 28.57, 88.56999999999999)
    5 = @(3.57, 100, 3.57, 92.86, 50, 89.29000000000001, 92.86, 21.43, 14.29, 21.43, 7.14, 21.43, 7.14, 10.71, 67.86)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $values.Length; $i++) {
        # Column B is index 2 (col offset starts at B = column 2)
        $col = $i + 2
        $ws.Cells.Item($row, $col).Value = $values[$i]
    }
}
